$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Agosto de 2020 a las 23:36'

$ws.Range("B4").Value = 6037197
$ws.Range("C4").Value = 36832
$ws.Range("D4").Value = 3335690
$ws.Range("E4").Value = 2517073
$ws.Range("G4").Value = 781
$ws.Range("H4").Value = 184434

$ws.Range("B5").Value = 3761391
$ws.Range("C5").Value = 39387
$ws.Range("E5").Value = 733894
$ws.Range("G5").Value = 893
$ws.Range("H5").Value = 118649

$ws.Range("B27").Value = 126672
$ws.Range("C27").Value = 255
$ws.Range("D27").Value = 112647
$ws.Range("E27").Value = 4926
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 9099

$ws.Range("B32").Value = 110403
$ws.Range("C32").Value = 2000
$ws.Range("D32").Value = 87726
$ws.Range("E32").Value = 21793
$ws.Range("G32").Value = 9
$ws.Range("H32").Value = 884

$ws.Range("B54").Value = 50756
$ws.Range("C54").Value = 363
$ws.Range("D54").Value = 47370
$ws.Range("E54").Value = 3198
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 188

$ws.Range("B79").Value = 17702
$ws.Range("C79").Value = 99
$ws.Range("D79").Value = 16139
$ws.Range("E79").Value = 1448
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 115

$ws.Range("B89").Value = 11601
$ws.Range("C89").Value = 225
$ws.Range("D89").Value = 10840
$ws.Range("E89").Value = 479

$ws.Range("A94").Value = 'Guinea'
$ws.Range("B94").Value = 9213
$ws.Range("C94").Value = 46
$ws.Range("D94").Value = 8180
$ws.Range("E94").Value = 975
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 58

$ws.Range("A95").Value = 'Croacia'
$ws.Range("B95").Value = 9192
$ws.Range("C95").Value = 304
$ws.Range("D95").Value = 6595
$ws.Range("E95").Value = 2420
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 177

$ws.Range("B104").Value = 6993
$ws.Range("C104").Value = 16
$ws.Range("D104").Value = 6381
$ws.Range("E104").Value = 454

$ws.Range("B107").Value = 5496
$ws.Range("C107").Value = 22
$ws.Range("D107").Value = 3121
$ws.Range("E107").Value = 2202

$ws.Range("A120").Value = 'Ruanda'
$ws.Range("B120").Value = 3672
$ws.Range("C120").Value = 47
$ws.Range("D120").Value = 1863
$ws.Range("E120").Value = 1794
$ws.Range("H120").Value = 15

$ws.Range("A121").Value = 'Mozambique'
$ws.Range("B121").Value = 3651
$ws.Range("C121").Value = 61
$ws.Range("D121").Value = 1968
$ws.Range("E121").Value = 1662
$ws.Range("H121").Value = 21

$ws.Range("A122").Value = 'Eslovaquia'
$ws.Range("B122").Value = 3626
$ws.Range("C122").Value = 90
$ws.Range("D122").Value = 2206
$ws.Range("E122").Value = 1387
$ws.Range("H122").Value = 33

$ws.Range("B163").Value = 1004
$ws.Range("C163").Value = 6
$ws.Range("D163").Value = 875
$ws.Range("E163").Value = 52

$ws.Range("B166").Value = 818
$ws.Range("C166").Value = 58
$ws.Range("E166").Value = 752
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 12

$ws.Range("D169").Value = 349
$ws.Range("E169").Value = 247

$ws.Range("B177").Value = 422
$ws.Range("C177").Value = 5
$ws.Range("E177").Value = 16

